# "update parameters for API"
#
# A new "Description" column is inserted between the "function" column and
# the "inputs"/"outputs" columns, documenting the REST endpoint (or a short
# note) for each function. The old "inputs" (C) and "outputs" (D) columns
# shift right to D and E. A couple of the (now) "outputs" cells also get
# extra/changed text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column C. Excel shifts the old inputs column
#    (C) to D and the old outputs column (D) to E, carrying over values,
#    styles and row heights automatically.
$ws.Columns("C:C").Insert()

# 2. Populate the new "Description" column (top-aligned, no wrap, like the
#    "function"/"inputs" columns) ...
$descriptions = @{
    "C1"  = "Description"
    "C3"  = "/segments/explore"
    "C5"  = "/segments/{id}/leaderboard"
}
foreach ($addr in $descriptions.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $descriptions[$addr]
    $cell.VerticalAlignment = -4160   # xlTop
    $cell.WrapText = $false
}

# ... except C6, which holds a longer note and wraps like the neighbouring
# inputs/outputs cells in that row.
$c6 = $ws.Range("C6")
$c6.Value = "compare elapsed time with segment length to compute  velocity, use minimum velocity theshold `noptional: use statistics to identify outliers (check distribution type!)"
$c6.VerticalAlignment = -4160
$c6.WrapText = $true

# The remaining new-column cells stay empty, but still top-aligned like the
# rest of the sheet.
foreach ($r in 2,4,7,8,9,10,11,12) {
    $cell = $ws.Range("C" + $r)
    $cell.VerticalAlignment = -4160
    $cell.WrapText = $false
}

# 3. Update the shifted "outputs" column (now E) with new/changed content.
$ws.Range("E4").Value = "list of segment IDs`nlist of segment coordinates`nsegment length in m"
$ws.Range("E6").Value = "number of entries for segment"

# 4. Fix up column widths for the new layout (closest values reachable
#    through the ColumnWidth property; Excel quantizes the saved width).
$ws.Columns("C:C").ColumnWidth = 42.5
$ws.Columns("D:D").ColumnWidth = 53.5
$ws.Columns("E:E").ColumnWidth = 27

# Rows 4 and 11 grow to fit the new/rewrapped text in the narrower columns.
$ws.Rows("4:4").RowHeight = 43.2
$ws.Rows("11:11").RowHeight = 43.2

# 5. Update the sheet view (zoom and selection).
$ws.Application.ActiveWindow.Zoom = 102
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C7").Select()
